# The sample watermark document originally carried its "Last update on
# ..." stamp in the page FOOTER (word/footer1.xml, right-aligned,
# footerReference in sectPr). The edit turns that stamp into a page
# HEADER instead (word/header1.xml, left-aligned, headerReference in
# sectPr) and refreshes the date it shows.

$d   = $word.ActiveDocument
$sec = $d.Sections.First

# ---------------------------------------------------------------------
# 1. Create/populate the primary header for the (only) section with the
#    content the footer used to hold, but left-aligned and with the
#    refreshed date. Section.Headers.Item(1) is wdHeaderFooterPrimary;
#    writing XML straight into its Range is the supported way to give a
#    header/footer story real content (Range.Text would also work, but
#    InsertXML lets us set the exact paragraph/run formatting in one
#    shot instead of relying on Find/Replace + Alignment afterwards).
$header = $sec.Headers.Item(1)
$headerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
             '<w:pPr><w:jc w:val="left"/></w:pPr>' +
             '<w:r>' +
             '<w:rPr><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' +
             '<w:t xml:space="preserve">Last update on 03/07/2022</w:t>' +
             '</w:r>' +
             '</w:p>'
$header.Range.InsertXML($headerXml)

# ---------------------------------------------------------------------
# 2. The old footer no longer carries the watermark stamp now that the
#    header owns it, so clear its text back out (the stamp shouldn't be
#    shown twice, once in the header and once in the footer).
$footer = $sec.Footers.Item(1)
$footer.Range.Text = ""
